$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.208689451217651
$ws.Range("B1").Value = 2.185441732406616
$ws.Range("C1").Value = 4.908904552459717
$ws.Range("D1").Value = 2.090900421142578
$ws.Range("E1").Value = 1.07283091545105
